$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-8 to append, columns A..L
$rows = @(
  ,@(520, "Leticia",    3,  "Cachorro", "Chihuahua",             "Branca",    "Mini",   510, "2024-06-11", "2024-06-11", "Históricos/520.txt", 0)
  ,@(2,   "Laila Said", 3,  "Cachorro", "Shar Pei",               "Caramelo",  "Médio",  2,   "2024-06-11", "2024-06-11", "Históricos/2.txt",   0)
  ,@(3,   "Lino",       5,  "Gato",     "Azul Russo",             "Caramelo",  "Grande", 3,   "2024-06-12", "2024-06-12", "Históricos/3.txt",   0)
  ,@(10,  "Coconut",    8,  "Cachorro", "Collie de Pelo Curto",   "Caramelo",  "Mini",   10,  "2024-06-12", "2024-06-12", "Históricos/10.txt",  0)
  ,@(11,  "Galvão",     13, "Gato",     "Egípcio Mau",            "Caramelo",  "Grande", 11,  "2024-06-12", "2024-06-12", "Históricos/11.txt",  0)
  ,@(12,  "Biscoito",   7,  "Gato",     "Exótico de Pelo Curto",  "Caramelo",  "Médio",  12,  "2024-06-12", "2024-06-12", "Históricos/12.txt",  0)
)

# Ensure date-looking text columns (I and J) are stored as text, not auto-converted dates
$ws.Range("I3:J8").NumberFormat = "@"

$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $r++
}
